$d = $word.ActiveDocument

# NOTE: the "Women's ice cream purchase frequency..." paragraph contains an
# apostrophe that must NOT be re-typed by Find/Replace, because Word's
# smart-quote AutoCorrect turns a straight apostrophe into a curly one
# whenever it falls inside a replaced span. So that one replacement below
# is scoped to the portion of the sentence strictly after the apostrophe.

$old0 = "Men: Men purchase ice cream both less than once a month and more than once a week, influenced by factors like deals, cravings, social context, and convenience, with brand preferences varying based on the occasion."
$new0 = "Men: Men purchase ice cream both less than once a month and more than once a week, influenced by factors like cravings, social context, and deals, with some prioritizing larger, economical tubs for frequent consumption."
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
Write-Host "Replace 0: $found0"
if (-not $found0) { throw "Replace 0 failed to find target text" }

$old1 = "variety, and brand loyalty, often viewing it as a simple pleasure or treat."
$new1 = "and variety, with some using it as a simple pleasure or reward."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Replace 1: $found1"
if (-not $found1) { throw "Replace 1 failed to find target text" }

$old2 = "Teens: Teens purchase ice cream from less than once a month to more than once a week, influenced by deals, social gatherings, unique flavors, and convenience, with brand preferences varying based on quality and personal taste."
$new2 = "Teens: Teens purchase ice cream from less than once a month to more than once a week, influenced by deals, social gatherings, and cravings, with some prioritizing unique flavors and smaller sizes for variety."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Replace 2: $found2"
if (-not $found2) { throw "Replace 2 failed to find target text" }

$old3 = "Men: Men often purchase ice cream based on promotions and value, sometimes opting for larger tubs at club stores or supermarkets. They also seek healthier options and are triggered by weather and post-meal cravings, occasionally indulging in convenient options like ice cream trucks."
$new3 = "Men: Men often purchase ice cream impulsively, driven by convenience and cravings, particularly when encountering ice cream carts or vending machines. They value quality ingredients and are open to trying new flavors, but are also influenced by promotions and price, especially when buying in bulk at club stores or supermarkets."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Replace 3: $found3"
if (-not $found3) { throw "Replace 3 failed to find target text" }

$old4 = "Women: Women appreciate variety in flavors and are influenced by convenience and price when purchasing ice cream. They often associate ice cream with emotional comfort and celebratory occasions, with brand loyalty varying based on flavor availability and quality."
$new4 = "Women: Women frequently associate ice cream with emotional comfort and indulgence, often purchasing it as a treat or reward. They appreciate readily available and affordable options from local grocery or convenience stores, and are drawn to brands with natural ingredients."
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replace 4: $found4"
if (-not $found4) { throw "Replace 4 failed to find target text" }

$old5 = "Teens: Teens are influenced by promotions and convenience, often purchasing ice cream from supermarkets and convenience stores. They enjoy trying new flavors and sharing with friends, with brand preferences varying based on taste and value."
$new5 = "Teens: Teens are influenced by multi-buy deals and promotions at supermarkets, and enjoy sharing and sampling different flavors with friends, especially when ordering online. They value convenience and are drawn to unique flavors, often purchasing ice cream from vending machines or corner stores as a quick treat."
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Replace 5: $found5"
if (-not $found5) { throw "Replace 5 failed to find target text" }
